$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DIP_Customer Information")

$ws.Range("B4").Value = "Test1"
$ws.Range("B5").Value = "Mining"
$ws.Range("B6").Value = "Nickel Mining"
$ws.Range("B8").Value = "Test2"
$ws.Range("B9").Value = "Test3"
$ws.Range("B10").Value = "Test4"
$ws.Range("B11").Value = "Test5"
$ws.Range("B12").Value = "Test6"
$ws.Range("B13").Value = "Test7"
$ws.Range("B16").Value = "Test8"
$ws.Range("C16").Value = "Test9"
$ws.Range("B17").Value = "Test10"
$ws.Range("C17").Value = "Test11"
$ws.Range("B18").Value = "Test12"
$ws.Range("C18").Value = "Test13"
$ws.Range("B19").Value = "Test14"
$ws.Range("C19").Value = "Test15"
$ws.Range("B20").Value = "Test16"
$ws.Range("C20").Value = "Test17"
$ws.Range("B21").Value = "Test18"
$ws.Range("C21").Value = "Test19"
$ws.Range("B22").Value = "Test20"
$ws.Range("C22").Value = "Test21"
$ws.Range("B23").Value = "Test22"
$ws.Range("C23").Value = "Test23"
$ws.Range("B24").Value = "Test24"
$ws.Range("B26").Value = "Test25"
$ws.Range("B28").Value = "Test26"
$ws.Range("B29").Value = "Test27"
$ws.Range("B30").Value = "Test28"
$ws.Range("B31").Value = "Test29"
$ws.Range("B32").Value = "Test30"
$ws.Range("B33").Value = "Test31"
$ws.Range("B35").Value = "Test32"
$ws.Range("C35").Value = "Test33"
$ws.Range("B36").Value = "Test34"
$ws.Range("C36").Value = "Test35"
$ws.Range("B37").Value = "Test36"
$ws.Range("C37").Value = "Test37"
$ws.Range("B38").Value = "Test38"
$ws.Range("C38").Value = "Test39"
$ws.Range("B39").Value = "Test40"
$ws.Range("C39").Value = "Test41"
$ws.Range("B40").Value = "Test42"
$ws.Range("C40").Value = "Test43"
$ws.Range("B41").Value = "Test44"
$ws.Range("C41").Value = "Test45"
$ws.Range("B42").Value = "Test46"
